$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Team (column B) values so the shared-string table is rebuilt in the target order
# and C column (PER-ish metric) values per the fixed calculation.
$ws.Cells.Item(2, 2).Value = "POR"
$ws.Cells.Item(2, 3).Value = 15.65333333333334
$ws.Cells.Item(3, 2).Value = "NJN"
$ws.Cells.Item(3, 3).Value = 11.425
$ws.Cells.Item(4, 2).Value = "CLE"
$ws.Cells.Item(4, 3).Value = 13.58
$ws.Cells.Item(5, 2).Value = "DAL"
$ws.Cells.Item(5, 3).Value = 14.10909090909091
$ws.Cells.Item(6, 2).Value = "MIA"
$ws.Cells.Item(6, 3).Value = 12.93333333333333
$ws.Cells.Item(7, 2).Value = "SEA"
$ws.Cells.Item(7, 3).Value = 14.65714285714286
$ws.Cells.Item(8, 2).Value = "ATL"
$ws.Cells.Item(8, 3).Value = 10.325
$ws.Cells.Item(9, 2).Value = "WAS"
$ws.Cells.Item(9, 3).Value = 12.64
$ws.Cells.Item(10, 2).Value = "MIL"
$ws.Cells.Item(10, 3).Value = 12.44666666666667
$ws.Cells.Item(11, 2).Value = "LAC"
$ws.Cells.Item(11, 3).Value = 13.43846153846154
$ws.Cells.Item(12, 2).Value = "VAN"
$ws.Cells.Item(12, 3).Value = 11.93076923076923
$ws.Cells.Item(13, 2).Value = "DET"
$ws.Cells.Item(13, 3).Value = 11.95
$ws.Cells.Item(14, 2).Value = "SAS"
$ws.Cells.Item(14, 3).Value = 14.18666666666667
$ws.Cells.Item(15, 2).Value = "ORL"
$ws.Cells.Item(15, 3).Value = 12.44666666666666
$ws.Cells.Item(16, 2).Value = "UTA"
$ws.Cells.Item(16, 3).Value = 13.42142857142857
$ws.Cells.Item(17, 2).Value = "HOU"
$ws.Cells.Item(17, 3).Value = 14.43846153846154
$ws.Cells.Item(18, 2).Value = "DEN"
$ws.Cells.Item(18, 3).Value = 12.60769230769231
$ws.Cells.Item(19, 2).Value = "LAL"
$ws.Cells.Item(19, 3).Value = 13.97333333333333
$ws.Cells.Item(20, 2).Value = "GSW"
$ws.Cells.Item(20, 3).Value = 10.50625
$ws.Cells.Item(21, 2).Value = "IND"
$ws.Cells.Item(21, 3).Value = 10.79333333333333
$ws.Cells.Item(22, 2).Value = "CHI"
$ws.Cells.Item(22, 3).Value = 10.50625
$ws.Cells.Item(23, 2).Value = "PHI"
$ws.Cells.Item(23, 3).Value = 13.53636363636364
$ws.Cells.Item(24, 2).Value = "CHH"
$ws.Cells.Item(24, 3).Value = 11.17333333333333
$ws.Cells.Item(25, 2).Value = "BOS"
$ws.Cells.Item(25, 3).Value = 12
$ws.Cells.Item(26, 2).Value = "TOR"
$ws.Cells.Item(26, 3).Value = 10.94
$ws.Cells.Item(27, 2).Value = "SAC"
$ws.Cells.Item(27, 3).Value = 13.65714285714286
$ws.Cells.Item(28, 2).Value = "PHO"
$ws.Cells.Item(28, 3).Value = 13.17272727272727
$ws.Cells.Item(29, 2).Value = "NYK"
$ws.Cells.Item(29, 3).Value = 12.15454545454545
$ws.Cells.Item(30, 2).Value = "MIN"
$ws.Cells.Item(30, 3).Value = 12.32142857142857
